$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.762.60'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '3.187.24'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.188.00'
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.51%  '
$ws.Range('D15').Value = '3.714.17'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '3.185.82'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').Value = '63.758.84'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('E19').Value = '  -1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '463.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('E22').Value = '  -1.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.71'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.92'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.82'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.58'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('E35').Value = '  -1.14%  '
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('D37').Value = '0.0₃0744'
$ws.Range('E37').Value = '  +6.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '51.70'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.70'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.15%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.115'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '400.08'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.68%  '
$ws.Range('D44').Value = '2.793.87'
$ws.Range('E44').Value = '  -7.02%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '35.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.75'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('E51').Value = '  +0.06%  '
